# Reversed DB query ordering on the index page changed which rows show up,
# and added new rows fetched from the (now reversed) result set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing rows whose values changed under the new ordering ---
$ws.Range("C2").Value = "CRISTAL TURISM SRL"
$ws.Range("D2").Value = 40
$ws.Range("D3").Value = 74

# --- New rows appended at the bottom (7-9) ---
$ws.Range("A7").Value = "maria.ioana.dicu@gmail.com"
$ws.Range("B7").Value = "Zona 3"
$ws.Range("C7").Value = "DEDEMAN"
$ws.Range("D7").Value = 344
$ws.Range("E7").Value = "'2025-08-01"
$ws.Range("E7").Style = "Normal"

$ws.Range("A8").Value = "maria.ioana.dicu@gmail.com"
$ws.Range("B8").Value = "Zona 4"
$ws.Range("C8").Value = "VIAMSO SRL"
$ws.Range("D8").Value = 100
$ws.Range("E8").Value = "'2025-08-12"
$ws.Range("E8").Style = "Normal"

$ws.Range("A9").Value = "maria.ioana.dicu@gmail.com"
$ws.Range("B9").Value = "Zona 2"
$ws.Range("C9").Value = "ARTHA STRUCTURE AG S.R.L"
$ws.Range("D9").Value = 567
$ws.Range("E9").Value = "'2025-08-12"
$ws.Range("E9").Style = "Normal"
